# Adding Append support to Excel.Write (#3558)
# -----------------------------------------------------------------
# This script reproduces, via Excel COM automation, the changes made to
# test/Table_Tests/data/TestSheet.xlsx:
#   1. Updates the selection on the "Another" worksheet to A2:D4 (active
#      cell D4).
#   2. Adds a new worksheet named "NoHeaders" at the end of the workbook,
#      containing the same data as "Another" (rows 2-4) but without a
#      header row, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1) Update selection on the "Another" sheet -------------------
$ws2 = $wb.Worksheets.Item("Another")
$ws2.Range("A2:D4").Select()

# --- 2) Add the new "NoHeaders" worksheet at the end ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "NoHeaders"

# Apply the date format to column D up front so the written dates reuse
# the workbook's existing "d-mmm" number format style instead of minting
# a brand-new (unused) one.
$newSheet.Range("D1:D3").NumberFormat = "d-mmm"

# Row 1: a / 1 / TRUE / 03-Jun-2022
$newSheet.Range("A1").Value = "a"
$newSheet.Range("B1").Value = 1
$newSheet.Range("C1").Value = $true
$newSheet.Range("D1").Value = (Get-Date -Year 2022 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 2: b / 2 / FALSE / 12-May-2022
$newSheet.Range("A2").Value = "b"
$newSheet.Range("B2").Value = 2
$newSheet.Range("C2").Value = $false
$newSheet.Range("D2").Value = (Get-Date -Year 2022 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)

# Row 3: c / 3 / FALSE / 15-Feb-2022
$newSheet.Range("A3").Value = "c"
$newSheet.Range("B3").Value = 3
$newSheet.Range("C3").Value = $false
$newSheet.Range("D3").Value = (Get-Date -Year 2022 -Month 2 -Day 15 -Hour 0 -Minute 0 -Second 0)

# Select the whole data range and make this the active sheet/tab, like
# the authored workbook (tabSelected moves from Sheet1 to NoHeaders).
$newSheet.Range("A1:D3").Select()
$newSheet.Activate()
